# Apply the "target create" bug fix: the Targets count (C16) increases from 2 to 10,
# which ripples through the dependent formulas in F18, C19, C20, and C21.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Targets (C16): 2 -> 10
$ws.Range("C16").Value = 10

# Make sure dependent formulas (F18, C19, C20, C21) are recalculated
$excel.Calculate()

# Leave the active selection on C17, matching the saved workbook state
$ws.Range("C17").Select()
